$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto market data.
# Values that are purely numeric-looking are prefixed with a leading
# apostrophe so Excel stores them as text (matching column D's existing text format),
# exactly as if a user had typed them into the cell.

$ws.Range("D2").Value = "63.914.85"
$ws.Range("E2").Value = "  -3.47%  "
$ws.Range("D3").Value = "3.129.64"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'610.14"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "'149.70"
$ws.Range("E6").Value = "  -4.43%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.130.11"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "  -3.67%  "
$ws.Range("E10").Value = "  -4.89%  "
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").Value = "'0.479"
$ws.Range("E12").Value = "  -5.21%  "
$ws.Range("D13").Value = "'0.0000260"
$ws.Range("E13").Value = "  -3.96%  "
$ws.Range("D14").Value = "'36.88"
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("D15").Value = "3.627.27"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").Value = "64.022.06"
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "3.116.49"
$ws.Range("E18").Value = "  -2.66%  "
$ws.Range("D19").Value = "'7.01"
$ws.Range("E19").Value = "  -4.64%  "
$ws.Range("D20").Value = "'484.11"
$ws.Range("E20").Value = "  -4.84%  "
$ws.Range("D21").Value = "'14.66"
$ws.Range("E21").Value = "  -4.38%  "
$ws.Range("D22").Value = "'0.713"
$ws.Range("E22").Value = "  -2.94%  "
$ws.Range("D23").Value = "'7.79"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("D24").Value = "'13.85"
$ws.Range("E24").Value = "  -5.56%  "
$ws.Range("D25").Value = "'84.27"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'2.94"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").Value = "'8.59"
$ws.Range("E28").Value = "  -5.42%  "
$ws.Range("D29").Value = "'0.126"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").Value = "'2.26"
$ws.Range("E30").Value = "  -4.35%  "
$ws.Range("D31").Value = "'7.05"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "'0.998"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").Value = "'2.71"
$ws.Range("E33").Value = "  -7.77%  "
$ws.Range("D34").Value = "'26.76"
$ws.Range("E34").Value = "  -5.37%  "
$ws.Range("D35").Value = "'1.11"
$ws.Range("E35").Value = "  -5.15%  "
$ws.Range("D36").Value = "'6.12"
$ws.Range("E36").Value = "  -6.00%  "
$ws.Range("D37").Value = "'54.53"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("E38").Value = "  +6.48%  "
$ws.Range("D39").Value = "0.0₃0748"
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("D40").Value = "'451.76"
$ws.Range("E40").Value = "  -10.03%  "
$ws.Range("E41").Value = "  -4.92%  "
$ws.Range("D42").Value = "'0.0401"
$ws.Range("E42").Value = "  -4.99%  "
$ws.Range("D43").Value = "'8.52"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D44").Value = "2.884.36"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").Value = "'0.273"
$ws.Range("E45").Value = "  -8.18%  "
$ws.Range("D46").Value = "'2.34"
$ws.Range("E46").Value = "  -4.13%  "
$ws.Range("D47").Value = "'26.79"
$ws.Range("E47").Value = "  -5.43%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").Value = "'2.32"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").Value = "'119.39"
$ws.Range("E51").Value = "  -2.56%  "

Write-Output "Updated crypto price/volume cells"
